# This workbook's rows 2-29 got reshuffled (weekly Fruta/Hortaliza refresh):
# for each destination row, the values of columns D (Fecha), J (Volumen),
# K (Precio minimo), L (Precio maximo), M (Precio promedio ponderado) and
# P (Precio $/Kg) are replaced by the values that used to sit in another
# source row (a pure permutation of the 28 data rows; row 21 is unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: destination row -> source row (values are read from the source
# row's *original* contents, so we snapshot everything first).
$map = @{
    2=25; 3=8; 4=12; 5=17; 6=20; 7=28; 8=4; 9=13; 10=3; 11=6;
    12=10; 13=22; 14=11; 15=24; 16=18; 17=16; 18=15; 19=7; 20=27; 21=21;
    22=26; 23=2; 24=19; 25=14; 26=23; 27=5; 28=29; 29=9
}

$cols = 4,10,11,12,13,16   # D, J, K, L, M, P

# Snapshot original values before writing anything, since several rows
# both give and receive values.
$snapshot = @{}
for ($r = 2; $r -le 29; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

foreach ($destRow in $map.Keys) {
    $srcRow = $map[$destRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value2 = $snapshot[$srcRow][$c]
    }
}
